$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new year headers (J4:N4), formatted like I4 (style index 3) ---
$ws.Range("I4").Copy() | Out-Null
$ws.Range("J4:N4").PasteSpecial(-4122) | Out-Null
$ws.Range("J4").Value = 2019
$ws.Range("K4").Value = 2020
$ws.Range("L4").Value = 2021
$ws.Range("M4").Value = 2022
$ws.Range("N4").Value = 2023

# --- Row 6: J6:N6, formatted like I6 (style index 7) ---
$ws.Range("I6").Copy() | Out-Null
$ws.Range("J6:N6").PasteSpecial(-4122) | Out-Null
$ws.Range("J6").Value = 81.1
$ws.Range("K6").Value = 85.8
$ws.Range("L6").Value = 78.1
$ws.Range("M6").Value = 72.2
$ws.Range("N6").Value = 75.7

# --- Row 7: J7:N7, formatted like I7 (style index 7) ---
$ws.Range("I7").Copy() | Out-Null
$ws.Range("J7:N7").PasteSpecial(-4122) | Out-Null
$ws.Range("J7").Value = 18.9
$ws.Range("K7").Value = 14.2
$ws.Range("L7").Value = 21.9
$ws.Range("M7").Value = 27.8
$ws.Range("N7").Value = 24.3

# --- Row 8: J8:N8, formatted like I8 (style index 7) ---
$ws.Range("I8").Copy() | Out-Null
$ws.Range("J8:N8").PasteSpecial(-4122) | Out-Null

# --- Row 9: J9:N9, formatted like I9 (style index 7) ---
$ws.Range("I9").Copy() | Out-Null
$ws.Range("J9:N9").PasteSpecial(-4122) | Out-Null
$ws.Range("J9").Value = 22.8
$ws.Range("K9").Value = 25.6
$ws.Range("L9").Value = 24.2
$ws.Range("M9").Value = 21.4
$ws.Range("N9").Value = 31.1

# --- Row 10: J10:N10, formatted like I10 (style index 7) ---
$ws.Range("I10").Copy() | Out-Null
$ws.Range("J10:N10").PasteSpecial(-4122) | Out-Null
$ws.Range("J10").Value = 77.2
$ws.Range("K10").Value = 74.4
$ws.Range("L10").Value = 75.8
$ws.Range("M10").Value = 78.6
$ws.Range("N10").Value = 68.9

# --- Row 11: J11:N11, formatted like I11 (style index 7) ---
$ws.Range("I11").Copy() | Out-Null
$ws.Range("J11:N11").PasteSpecial(-4122) | Out-Null

# --- Row 12: J12:N12, formatted like I12 (style index 7) ---
$ws.Range("I12").Copy() | Out-Null
$ws.Range("J12:N12").PasteSpecial(-4122) | Out-Null
$ws.Range("J12").Value = 84.4
$ws.Range("K12").Value = 72.7
$ws.Range("L12").Value = 73.3
$ws.Range("M12").Value = 72.8
$ws.Range("N12").Value = 76.7

# --- Row 13: J13:N13, formatted like I13 (style index 7) ---
$ws.Range("I13").Copy() | Out-Null
$ws.Range("J13:N13").PasteSpecial(-4122) | Out-Null
$ws.Range("J13").Value = 15.6
$ws.Range("K13").Value = 27.3
$ws.Range("L13").Value = 26.7
$ws.Range("M13").Value = 27.2
$ws.Range("N13").Value = 23.3

# --- Row 14: J14:N14, formatted like I14 (style index 7) ---
$ws.Range("I14").Copy() | Out-Null
$ws.Range("J14:N14").PasteSpecial(-4122) | Out-Null

# --- Row 15: J15:N15, formatted like I15 (style index 7) ---
$ws.Range("I15").Copy() | Out-Null
$ws.Range("J15:N15").PasteSpecial(-4122) | Out-Null
$ws.Range("J15").Value = 90.3
$ws.Range("K15").Value = 93.4
$ws.Range("L15").Value = 90.5
$ws.Range("M15").Value = 87.8
$ws.Range("N15").Value = 89

# --- Row 16: J16:N16, formatted like I16 (style index 7) ---
$ws.Range("I16").Copy() | Out-Null
$ws.Range("J16:N16").PasteSpecial(-4122) | Out-Null
$ws.Range("J16").Value = 9.7
$ws.Range("K16").Value = 6.6
$ws.Range("L16").Value = 9.5
$ws.Range("M16").Value = 12.2
$ws.Range("N16").Value = 11

# --- Row 17: J17:N17, formatted like I17 (style index 7) ---
$ws.Range("I17").Copy() | Out-Null
$ws.Range("J17:N17").PasteSpecial(-4122) | Out-Null

# --- Row 18: J18:N18, formatted like I18 (style index 7) ---
$ws.Range("I18").Copy() | Out-Null
$ws.Range("J18:N18").PasteSpecial(-4122) | Out-Null
$ws.Range("J18").Value = 60.2
$ws.Range("K18").Value = 66
$ws.Range("L18").Value = 59.3
$ws.Range("M18").Value = 44.9
$ws.Range("N18").Value = 48.3

# --- Row 19: J19:N19, formatted like I19 (style index 8) ---
$ws.Range("I19").Copy() | Out-Null
$ws.Range("J19:N19").PasteSpecial(-4122) | Out-Null
$ws.Range("J19").Value = 39.8
$ws.Range("K19").Value = 34
$ws.Range("L19").Value = 40.7
$ws.Range("M19").Value = 55.1
$ws.Range("N19").Value = 51.7

# --- Row 20: new row height + new font style (size 8 Times New Roman) for A20:C20 ---
$ws.Rows.Item(20).RowHeight = 15.75
$ws.Range("A20:C20").Font.Name = "Times New Roman"
$ws.Range("A20:C20").Font.Size = 8

# --- Clear clipboard ---
$excel.CutCopyMode = $false
